# Insert a new weekly price record as row 43, shifting the existing
# rows 43-58 down to 44-59 (dimension grows from A1:R58 to A1:R59).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43:43").Insert()

$ws.Range("A43").Value = 6
$ws.Range("B43").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44841
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112035
$ws.Range("G43").Value = "Bruselas (repollito)"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 580
$ws.Range("K43").Value = 17000
$ws.Range("L43").Value = 18000
$ws.Range("M43").Value = 17448
$ws.Range("N43").Value = "$/malla 15 kilos"
$ws.Range("O43").Value = "Provincia de Quillota"
$ws.Range("P43").Value = 1163
$ws.Range("Q43").Value = 15
$ws.Range("R43").Value = "Hortaliza"
